# Update countries & provincias Spain
# - Swap ranking order of Brasil/Alemania (rows 10-11) and Argentina/Marruecos
#   (rows 55-56) to reflect updated case counts, refresh the "datos
#   actualizados" timestamp, and bump the various per-country case/death
#   counters to their newer values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 00:35"
$ws.Range("B4").Value = 1406807
$ws.Range("C4").Value = 20973
$ws.Range("D4").Value = 280509
$ws.Range("E4").Value = 1043014
$ws.Range("G4").Value = 1489
$ws.Range("H4").Value = 83284
$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = 177589
$ws.Range("C10").Value = 8446
$ws.Range("D10").Value = 72597
$ws.Range("E10").Value = 92592
$ws.Range("F10").Value = 8318
$ws.Range("G10").Value = 775
$ws.Range("H10").Value = 12400
$ws.Range("A11").Value = "Alemania"
$ws.Range("B11").Value = 173171
$ws.Range("C11").Value = 595
$ws.Range("D11").Value = 147200
$ws.Range("E11").Value = 18233
$ws.Range("F11").Value = 1539
$ws.Range("G11").Value = 77
$ws.Range("H11").Value = 7738
$ws.Range("B17").Value = 71157
$ws.Range("C17").Value = 1176
$ws.Range("D17").Value = 34042
$ws.Range("E17").Value = 31946
$ws.Range("G17").Value = 176
$ws.Range("H17").Value = 5169
$ws.Range("B52").Value = 8157
$ws.Range("C52").Value = 25
$ws.Range("E52").Value = 7897
$ws.Range("A55").Value = "Argentina"
$ws.Range("B55").Value = 6563
$ws.Range("C55").Value = 285
$ws.Range("D55").Value = 1862
$ws.Range("E55").Value = 4382
$ws.Range("F55").Value = 170
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 319
$ws.Range("A56").Value = "Marruecos"
$ws.Range("B56").Value = 6418
$ws.Range("C56").Value = 137
$ws.Range("D56").Value = 2991
$ws.Range("E56").Value = 3239
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 188
$ws.Range("D135").Value = 272
$ws.Range("E135").Value = 36
